# Updated cryptos list on Sun Sep  1 03:48:59 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.557.58"
$ws.Range("E2").Value = "  -1.15%  "
$ws.Range("D3").Value = "2.487.24"
$ws.Range("E3").Value = "  -1.51%  "
$ws.Range("E4").Value = "  +0.34%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "527.55"
$ws.Range("E5").Value = "  -1.91%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.01"
$ws.Range("E6").Value = "  -3.29%  "
$ws.Range("E7").Value = "  +0.24%  "
$ws.Range("E8").Value = "  -0.85%  "
$ws.Range("E9").Value = "  -0.92%  "
$ws.Range("E10").Value = "  -2.10%  "
$ws.Range("E11").Value = "  +0.53%  "
$ws.Range("E12").Value = "  -0.94%  "
$ws.Range("D13").Value = "2.928.13"
$ws.Range("E13").Value = "  -1.58%  "
$ws.Range("D14").Value = "58.459.74"
$ws.Range("E14").Value = "  -1.09%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "22.49"
$ws.Range("E15").Value = "  -3.21%  "
$ws.Range("E16").Value = "  -1.89%  "
$ws.Range("D17").Value = "2.487.84"
$ws.Range("E17").Value = "  -1.14%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.95"
$ws.Range("E18").Value = "  -1.46%  "
$ws.Range("E19").Value = "  -1.88%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "321.97"
$ws.Range("E20").Value = "  -1.37%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.999"
$ws.Range("E21").Value = "  +0.04%  "
$ws.Range("E22").Value = "  -1.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "64.44"
$ws.Range("E23").Value = "  -1.90%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.414"
$ws.Range("E24").Value = "  -2.24%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("E25").Value = "  +0.12%  "
$ws.Range("E26").Value = "  -2.44%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.48"
$ws.Range("E27").Value = "  -2.47%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.43"
$ws.Range("E29").Value = "  -5.05%  "
$ws.Range("E30").Value = "  -3.56%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "165.14"
$ws.Range("E31").Value = "  -2.60%  "
$ws.Range("E32").Value = "  -5.15%  "
$ws.Range("E33").Value = "  +0.10%  "
$ws.Range("E34").Value = "  -0.07%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "18.26"
$ws.Range("E35").Value = "  -1.59%  "
$ws.Range("E36").Value = "  -8.71%  "
$ws.Range("E38").Value = "  -4.01%  "
$ws.Range("E39").Value = "  -2.79%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.53"
$ws.Range("E40").Value = "  -3.17%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "276.75"
$ws.Range("E41").Value = "  -2.66%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.96"
$ws.Range("E42").Value = "  -5.85%  "
$ws.Range("E43").Value = "  -1.46%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "127.44"
$ws.Range("E44").Value = "  -2.33%  "
$ws.Range("E45").Value = "  -1.73%  "
$ws.Range("E46").Value = "  -3.03%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "17.26"
$ws.Range("E48").Value = "  -1.73%  "
$ws.Range("D49").Value = "1.742.25"
$ws.Range("E49").Value = "  -1.41%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.975"
$ws.Range("E50").Value = "  -1.38%  "
$ws.Range("E51").Value = "  -1.75%  "
